# Auto-generated edit script applying scheduled Chocobo_Profits market-data refresh
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across several
# crafting-class sheets with refreshed Universalis price data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 459
$ws.Range("J28").Value = 464.5
$ws.Range("L28").Value = 464.5
$ws.Range("N28").Value = -1434.5
$ws.Range("H107").Value = 1658
$ws.Range("J107").Value = 1037.25
$ws.Range("L107").Value = 1037.25
$ws.Range("N107").Value = -4877.25
$ws.Range("H131").Value = 3611.7896
$ws.Range("I131").Value = 1766.1111
$ws.Range("J131").Value = 5272.9
$ws.Range("K131").Value = 5298.3333
$ws.Range("L131").Value = 15818.7
$ws.Range("M131").Value = -258.3333000000002
$ws.Range("N131").Value = -25898.7
$ws.Range("H132").Value = 13838638
$ws.Range("I132").Value = 17860222
$ws.Range("J132").Value = 591066.4399999999
$ws.Range("K132").Value = 53580666
$ws.Range("L132").Value = 1773199.32
$ws.Range("M132").Value = -53578136
$ws.Range("N132").Value = -1778259.32
$ws.Range("H135").Value = 670.53845
$ws.Range("I135").Value = 293.5
$ws.Range("J135").Value = 1518.875
$ws.Range("K135").Value = 2641.5
$ws.Range("L135").Value = 13669.875
$ws.Range("M135").Value = -106.5
$ws.Range("N135").Value = -18739.875
$ws.Range("H137").Value = 1829.431
$ws.Range("I137").Value = 938.65216
$ws.Range("J137").Value = 5244.0835
$ws.Range("K137").Value = 2815.95648
$ws.Range("L137").Value = 15732.2505
$ws.Range("M137").Value = -265.9564799999998
$ws.Range("N137").Value = -20832.2505
$ws.Range("H138").Value = 2413.75
$ws.Range("I138").Value = 755.4054
$ws.Range("J138").Value = 3387.6985
$ws.Range("K138").Value = 2266.2162
$ws.Range("L138").Value = 10163.0955
$ws.Range("M138").Value = 2873.7838
$ws.Range("N138").Value = -20443.0955
$ws.Range("H141").Value = 5007.5093
$ws.Range("I141").Value = 5563.8184
$ws.Range("J141").Value = 2287.7778
$ws.Range("K141").Value = 16691.4552
$ws.Range("L141").Value = 6863.3334
$ws.Range("M141").Value = -11511.4552
$ws.Range("N141").Value = -17223.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4341.7144
$ws.Range("I32").Value = 3074.403
$ws.Range("K32").Value = 3074.403
$ws.Range("M32").Value = -2787.403
$ws.Range("H61").Value = 1193.7916
$ws.Range("I61").Value = 1015.7857
$ws.Range("J61").Value = 1443
$ws.Range("K61").Value = 1015.7857
$ws.Range("L61").Value = 1443
$ws.Range("M61").Value = -803.7857
$ws.Range("N61").Value = -1867
$ws.Range("H74").Value = 2526.52
$ws.Range("I74").Value = 2528.6924
$ws.Range("K74").Value = 2528.6924
$ws.Range("M74").Value = -1654.6924
$ws.Range("H77").Value = 2526.52
$ws.Range("I77").Value = 2528.6924
$ws.Range("K77").Value = 12643.462
$ws.Range("M77").Value = -8275.462
$ws.Range("H132").Value = 1711.8928
$ws.Range("I132").Value = 1073.909
$ws.Range("J132").Value = 4051.1667
$ws.Range("K132").Value = 3221.727
$ws.Range("L132").Value = 12153.5001
$ws.Range("M132").Value = -691.7270000000003
$ws.Range("N132").Value = -17213.5001
$ws.Range("H136").Value = 1193.7916
$ws.Range("I136").Value = 1015.7857
$ws.Range("J136").Value = 1443
$ws.Range("K136").Value = 3047.3571
$ws.Range("L136").Value = 4329
$ws.Range("M136").Value = -497.3571000000002
$ws.Range("N136").Value = -9429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 24319.5
$ws.Range("I82").Value = 4571
$ws.Range("J82").Value = 29705.455
$ws.Range("K82").Value = 4571
$ws.Range("L82").Value = 29705.455
$ws.Range("M82").Value = -4188
$ws.Range("N82").Value = -30471.455
$ws.Range("H85").Value = 24319.5
$ws.Range("I85").Value = 4571
$ws.Range("J85").Value = 29705.455
$ws.Range("K85").Value = 4571
$ws.Range("L85").Value = 29705.455
$ws.Range("M85").Value = -3245
$ws.Range("N85").Value = -32357.455
$ws.Range("H94").Value = 934.75
$ws.Range("I94").Value = 1106.25
$ws.Range("J94").Value = 591.75
$ws.Range("K94").Value = 1106.25
$ws.Range("L94").Value = 591.75
$ws.Range("M94").Value = -655.25
$ws.Range("N94").Value = -1493.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 22223682
$ws.Range("J16").Value = 1300
$ws.Range("L16").Value = 1300
$ws.Range("N16").Value = -1874
$ws.Range("H31").Value = 9617649
$ws.Range("I31").Value = 1280.75
$ws.Range("J31").Value = 31254476
$ws.Range("K31").Value = 1280.75
$ws.Range("L31").Value = 31254476
$ws.Range("M31").Value = -985.75
$ws.Range("N31").Value = -31255066
$ws.Range("H34").Value = 9617649
$ws.Range("I34").Value = 1280.75
$ws.Range("J34").Value = 31254476
$ws.Range("K34").Value = 1280.75
$ws.Range("L34").Value = 31254476
$ws.Range("M34").Value = -1078.75
$ws.Range("N34").Value = -31254880
$ws.Range("H58").Value = 1813.8379
$ws.Range("I58").Value = 1553.6716
$ws.Range("J58").Value = 4304
$ws.Range("K58").Value = 1553.6716
$ws.Range("L58").Value = 4304
$ws.Range("M58").Value = -1350.6716
$ws.Range("N58").Value = -4710
$ws.Range("H107").Value = 813.5333000000001
$ws.Range("I107").Value = 524
$ws.Range("J107").Value = 1392.6
$ws.Range("K107").Value = 524
$ws.Range("L107").Value = 1392.6
$ws.Range("M107").Value = 1396
$ws.Range("N107").Value = -5232.6
$ws.Range("H113").Value = 22223682
$ws.Range("J113").Value = 1300
$ws.Range("L113").Value = 1300
$ws.Range("N113").Value = -5640
$ws.Range("H122").Value = 2904.7273
$ws.Range("I122").Value = 1729.3334
$ws.Range("J122").Value = 3345.5
$ws.Range("K122").Value = 5188.0002
$ws.Range("L122").Value = 10036.5
$ws.Range("M122").Value = -2738.0002
$ws.Range("N122").Value = -14936.5
$ws.Range("H134").Value = 3340.6226
$ws.Range("I134").Value = 3544.611
$ws.Range("J134").Value = 2908.647
$ws.Range("K134").Value = 10633.833
$ws.Range("L134").Value = 8725.940999999999
$ws.Range("M134").Value = -8098.832999999999
$ws.Range("N134").Value = -13795.941
$ws.Range("H136").Value = 1813.8379
$ws.Range("I136").Value = 1553.6716
$ws.Range("J136").Value = 4304
$ws.Range("K136").Value = 4661.0148
$ws.Range("L136").Value = 12912
$ws.Range("M136").Value = -2111.0148
$ws.Range("N136").Value = -18012

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2551108
$ws.Range("I2").Value = 90
$ws.Range("J2").Value = 2747340
$ws.Range("K2").Value = 540
$ws.Range("L2").Value = 16484040
$ws.Range("M2").Value = -427
$ws.Range("N2").Value = -16484266
$ws.Range("H113").Value = 876.72
$ws.Range("I113").Value = 731.6875
$ws.Range("J113").Value = 1134.5555
$ws.Range("K113").Value = 2195.0625
$ws.Range("L113").Value = 3403.6665
$ws.Range("M113").Value = -25.0625
$ws.Range("N113").Value = -7743.666499999999
$ws.Range("H137").Value = 2269.2222
$ws.Range("J137").Value = 3170.5862
$ws.Range("L137").Value = 9511.758600000001
$ws.Range("N137").Value = -19711.7586

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M46").ClearContents()
$ws.Range("H46").Value = 32304
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 32304
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 32304
$ws.Range("N46").Value = -32616
$ws.Range("H126").Value = 2536.04
$ws.Range("I126").Value = 2536.04
$ws.Range("K126").Value = 7608.12
$ws.Range("M126").Value = -5138.12
$ws.Range("H132").Value = 2542.3948
$ws.Range("I132").Value = 1516.16
$ws.Range("K132").Value = 4548.48
$ws.Range("M132").Value = -2018.48

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 6285.7144
$ws.Range("J14").Value = 6285.7144
$ws.Range("L14").Value = 6285.7144
$ws.Range("N14").Value = -6621.7144
$ws.Range("H132").Value = 11907239
$ws.Range("I132").Value = 1843.8889
$ws.Range("J132").Value = 33336948
$ws.Range("K132").Value = 5531.6667
$ws.Range("L132").Value = 100010844
$ws.Range("M132").Value = -3001.6667
$ws.Range("N132").Value = -100015904
$ws.Range("H136").Value = 2080.7917
$ws.Range("I136").Value = 633.0540999999999
$ws.Range("J136").Value = 6950.4546
$ws.Range("K136").Value = 1899.1623
$ws.Range("L136").Value = 20851.3638
$ws.Range("M136").Value = 650.8377
$ws.Range("N136").Value = -25951.3638
